$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-53 (values are recomputed by the naive forecaster; row 2 is newly inserted)
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 5.896808312953783
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 7.318442086255605
$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 6.056254825277896
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 4.566338461218011
$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 2008
$ws.Cells.Item(4, 3).Value = 7.441962824572235
$ws.Cells.Item(4, 4).Value = 2009
$ws.Cells.Item(4, 5).Value = 8.262942840582955
$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = 8.604123301398037
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = 8.260999835306727
$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2009
$ws.Cells.Item(6, 3).Value = 6.277541464866987
$ws.Cells.Item(6, 4).Value = 2010
$ws.Cells.Item(6, 5).Value = 7.915558093865038
$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 5.436647924209592
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = 4.862860110364875
$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 2010
$ws.Cells.Item(8, 3).Value = 6.535114773304773
$ws.Cells.Item(8, 4).Value = 2011
$ws.Cells.Item(8, 5).Value = 6.402044794134309
$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 6.334380382529425
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = 5.829578861489648
$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 2011
$ws.Cells.Item(10, 3).Value = 5.12051970717502
$ws.Cells.Item(10, 4).Value = 2012
$ws.Cells.Item(10, 5).Value = 5.160655543917292
$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 3.88993859232436
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 3.4300351921007
$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 2012
$ws.Cells.Item(12, 3).Value = 3.65682115264816
$ws.Cells.Item(12, 4).Value = 2013
$ws.Cells.Item(12, 5).Value = 3.62442670409151
$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = 2.513767348245044
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 2.479992751939486
$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 2013
$ws.Cells.Item(14, 3).Value = 2.943878639034381
$ws.Cells.Item(14, 4).Value = 2014
$ws.Cells.Item(14, 5).Value = 3.771815305047821
$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 2.723916849952834
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 2.551173534479334
$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 2014
$ws.Cells.Item(16, 3).Value = 1.172679597477866
$ws.Cells.Item(16, 4).Value = 2015
$ws.Cells.Item(16, 5).Value = 1.204704113773114
$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = 1.326505206336948
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 1.211929054838756
$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 2.961845079861303
$ws.Cells.Item(18, 4).Value = 2016
$ws.Cells.Item(18, 5).Value = 2.562199564969392
$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 2.321003614014883
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 2.610201636760778
$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 2016
$ws.Cells.Item(20, 3).Value = 2.508469427909898
$ws.Cells.Item(20, 4).Value = 2017
$ws.Cells.Item(20, 5).Value = 2.632055757778851
$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 2.468891199411116
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 2.593292206016984
$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = 3.523703831572056
$ws.Cells.Item(22, 4).Value = 2018
$ws.Cells.Item(22, 5).Value = 3.331068508781954
$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 2.442951431721241
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 2.874490643722805
$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 3.133596157287766
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 3.419422858788335
$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 1.273091634877033
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 1.947147525128035
$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 2018
$ws.Cells.Item(26, 3).Value = 1.178605266817589
$ws.Cells.Item(26, 4).Value = 2019
$ws.Cells.Item(26, 5).Value = 1.467147844249106
$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 1.388614840712377
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 1.965204165904111
$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 0.08486825492834971
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 0.9311475558545057
$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 2.983246785467752
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 2.844177684788551
$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = 3.047037961814492
$ws.Cells.Item(30, 4).Value = 2020
$ws.Cells.Item(30, 5).Value = 2.757054249287738
$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 3.192446966735796
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 2.739953120486138
$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = 2.405224065057476
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = 2.152035263856344
$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -0.4891791466461126
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = 0.6596954352570572
$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = -0.2228847697281378
$ws.Cells.Item(34, 4).Value = 2021
$ws.Cells.Item(34, 5).Value = 2.031388171904314
$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = 2.00987424531256
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = 1.446144849547082
$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = -0.678826357714013
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = 0.02414656897629097
$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = -0.7941560676977599
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = 1.553668892952742
$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 2021
$ws.Cells.Item(38, 3).Value = -1.165854108406617
$ws.Cells.Item(38, 4).Value = 2022
$ws.Cells.Item(38, 5).Value = 2.209187654930855
$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 3.122586570118835
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = 1.070586793467498
$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 0.7010162698181555
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = -0.6730402944081559
$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 2.34069710769782
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = 1.423433561452137
$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 2022
$ws.Cells.Item(42, 3).Value = 2.501311189006916
$ws.Cells.Item(42, 4).Value = 2023
$ws.Cells.Item(42, 5).Value = 2.819423640661167
$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = 2.639029809913129
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = 1.61912427117199
$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = 2.158153176293576
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = 1.298013848993262
$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = 0.9259311313598806
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = 1.321448002249337
$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 2023
$ws.Cells.Item(46, 3).Value = 0.6753076481029074
$ws.Cells.Item(46, 4).Value = 2024
$ws.Cells.Item(46, 5).Value = 0.6583865490665364
$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = 1.29262043481877
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = 1.577641784671369
$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = -0.1645072558042915
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 0.1619933518385297
$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = 1.935025917091848
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = 1.814429727910216
$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 2024
$ws.Cells.Item(50, 3).Value = 2.039329803030121
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = 2.453497774623137
$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = 1.562469275846001
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = 1.276308453105246
$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = 3.131832690451031
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = 2.199676451050503
$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 2025
$ws.Cells.Item(53, 3).Value = 2.43119486791763
$ws.Cells.Item(53, 4).Value = 2026
$ws.Cells.Item(53, 5).Value = 2.795512766421537

# Row 53 is brand new; copy the date-column number format/style from row 52 so A53 matches the others
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0
